$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (collector_name) for the new "village" field.
$ws.Columns.Item(5).Insert()

# The newly inserted column should keep the same width as the column to its
# left (admin_level_2, column D), matching Excel's default insert behaviour.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Set the new header label.
$ws.Range("E1").Value = "village"
